# THE BD IS REAL
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New H-column values (employee-entreprise association ids) for rows 2..50
$values = @{
    2  = 25
    3  = 7
    4  = 30
    5  = 48
    6  = 28
    7  = 16
    8  = 22
    9  = 41
    10 = 24
    11 = 2
    12 = 22
    13 = 39
    14 = 46
    16 = 2
    17 = 10
    18 = 33
    19 = 43
    20 = 15
    21 = 10
    22 = 34
    23 = 49
    24 = 37
    25 = 8
    26 = 46
    27 = 11
    28 = 44
    29 = 29
    30 = 6
    31 = 39
    32 = 17
    33 = 29
    34 = 21
    35 = 34
    36 = 1
    37 = 39
    38 = 14
    39 = 48
    40 = 25
    41 = 15
    42 = 14
    43 = 41
    44 = 44
    45 = 45
    46 = 24
    47 = 21
    49 = 22
    50 = 11
}

foreach ($row in $values.Keys) {
    $ws.Range("H$row").Value = $values[$row]
}

# Update the selection/view state saved with the worksheet
$ws.Range("K6").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 3
